$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph that declares "id_produto INT PRIMARY KEY
# AUTO_INCREMENT," used a paragraph-level first-line indent (w:ind
# firstLine="708") to visually offset the line under "CREATE TABLE
# tb_produto(". The revision instead removes that indent and inserts a
# literal Tab character run in front of the text, matching how the sibling
# "INSERT INTO" block already indents its continuation line with an
# explicit <w:tab/> run.
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("id_produto INT PRIMARY KEY AUTO_INCREMENT,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $p1 = $r1.Paragraphs.Item(1)
    $p1.Format.FirstLineIndent = 0
    $insertionPoint = $d.Range($r1.Start, $r1.Start)
    $insertionPoint.InsertBefore("`t")
}

# ---------------------------------------------------------------------
# Change 2: question 30's SELECT now asks for marca + a total-count alias
# instead of selecting the (unused) estoque_disponivel column.
# Use a direct Range.Text assignment (after locating the text with Find)
# rather than Find.Execute's ReplaceWith, because ReplaceWith silently
# "smart quotes"/autocorrects straight backticks/apostrophes.
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("SELECT COUNT(nome), marca, estoque_disponivel FROM", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Text = 'SELECT marca, COUNT(nome) AS `total de produto` FROM'
}

# ---------------------------------------------------------------------
# Change 3: the HAVING clause now filters on the aggregated COUNT(nome)
# instead of the removed estoque_disponivel column.
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("HAVING estoque_disponivel > '5';", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Text = "HAVING COUNT(nome) > '5';"
}
